$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.837.02'
$ws.Range('E2').Value = '  -4.53%  '

$ws.Range('D3').Value = '3.138.16'
$ws.Range('E3').Value = '  -6.46%  '

$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '515.66'
$ws.Range('E5').Value = '  -1.88%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '167.17'
$ws.Range('E6').Value = '  -8.27%  '

$ws.Range('E7').Value = '  -2.05%  '

$ws.Range('E8').Value = '  +0.05%  '

$ws.Range('D9').Value = '3.139.62'
$ws.Range('E9').Value = '  -6.20%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.590'
$ws.Range('E10').Value = '  -3.88%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '51.34'
$ws.Range('E11').Value = '  -8.93%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.127'
$ws.Range('E12').Value = '  -2.37%  '

$ws.Range('E13').Value = '  -2.28%  '

$ws.Range('E14').Value = '  -2.79%  '

$ws.Range('D15').Value = '3.638.36'
$ws.Range('E15').Value = '  -6.54%  '

$ws.Range('E16').Value = '  -5.67%  '

$ws.Range('D17').Value = '3.135.93'
$ws.Range('E17').Value = '  -6.70%  '

$ws.Range('D18').Value = '61.713.06'
$ws.Range('E18').Value = '  -4.37%  '

$ws.Range('E19').Value = '  -2.21%  '

$ws.Range('B20').Value = 'Polygon'
$ws.Range('C20').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.955'
$ws.Range('E20').Value = '  +0.35%  '

$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.75'
$ws.Range('E21').Value = '  -1.17%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '357.14'
$ws.Range('E22').Value = '  -2.95%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.03'
$ws.Range('E23').Value = '  +3.68%  '

$ws.Range('B24').Value = 'PancakeSwap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.65'
$ws.Range('E24').Value = '  -0.85%  '

$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '79.37'
$ws.Range('E25').Value = '  -1.13%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.87'
$ws.Range('E26').Value = '  +5.21%  '

$ws.Range('E27').Value = '  +4.08%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.56'
$ws.Range('E28').Value = '  -1.79%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '11.03'
$ws.Range('E29').Value = '  -1.17%  '

$ws.Range('E30').Value = '  -4.46%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '631.37'
$ws.Range('E31').Value = '  -5.06%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '27.65'
$ws.Range('E32').Value = '  -4.77%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.27'
$ws.Range('E33').Value = '  -6.04%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.12'
$ws.Range('E34').Value = '  +1.06%  '

$ws.Range('E35').Value = '  -0.37%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '56.01'
$ws.Range('E36').Value = '  -7.15%  '

$ws.Range('E37').Value = '  +0.03%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '36.01'
$ws.Range('E38').Value = '  +0.26%  '

$ws.Range('E39').Value = '  -1.71%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.997'

$ws.Range('D41').Value = '0.0₃0675'
$ws.Range('E41').Value = '  +10.78%  '

$ws.Range('E42').Value = '  -2.39%  '

$ws.Range('D43').Value = '2.859.99'
$ws.Range('E43').Value = '  +2.46%  '

$ws.Range('E44').Value = '  +8.54%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.85'
$ws.Range('E45').Value = '  +12.06%  '

$ws.Range('E46').Value = '  -0.25%  '

$ws.Range('E47').Value = '  -0.57%  '

$ws.Range('B48').Value = 'ApeXProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.95'
$ws.Range('E48').Value = '  +5.34%  '

$ws.Range('B49').Value = 'ThetaToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.49'
$ws.Range('E49').Value = '  -8.48%  '

$ws.Range('E50').Value = '  -2.29%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '131.12'
$ws.Range('E51').Value = '  -3.09%  '
